# BOM.xlsx update: "Added zips for gerber files"
# Adds three new connector header parts (rows 24-26), a note in row 20 (C20),
# and moves the Total row down to row 32 (summing G3:G26), reflecting a
# Currency cell-style for the whole G column instead of the ad-hoc format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 24: 5 position female header
# ---------------------------------------------------------------------
$ws.Range("B21:G21").Copy()
$ws.Range("B24:G24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B24").Value = "PPPC051LFBN-RC"
$ws.Range("C24").Value = "5 position female header"
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 0.68
$ws.Range("F24").Value = "https://www.digikey.ca/product-detail/en/sullins-connector-solutions/PPPC051LFBN-RC/S7038-ND/810177"
$ws.Range("G24").Formula = "=E24*D24"

# ---------------------------------------------------------------------
# Row 25: 10 position female header
# ---------------------------------------------------------------------
$ws.Range("B21:G21").Copy()
$ws.Range("B25:G25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B25").Value = "PPTC101LFBN-RC"
$ws.Range("C25").Value = "10 position female header"
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 0.9
$ws.Range("F25").Value = "https://www.digikey.ca/product-detail/en/sullins-connector-solutions/PPTC101LFBN-RC/S7008-ND/810149"
$ws.Range("G25").Formula = "=E25*D25"

# ---------------------------------------------------------------------
# Row 26: 10 position xbee header (this used to be the Total row)
# ---------------------------------------------------------------------
$ws.Range("B21:G21").Copy()
$ws.Range("B26:G26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B26").Value = "NPPN101BFCN-RC"
$ws.Range("C26").Value = "10 position xbee header"
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 1.36
$ws.Range("F26").Value = "https://www.digikey.ca/product-detail/en/sullins-connector-solutions/NPPN101BFCN-RC/S5751-10-ND/804812"
$ws.Range("G26").Formula = "=E26*D26"

# ---------------------------------------------------------------------
# G column: switch from the ad-hoc currency format over to the Currency
# cell-style (same look, matches column E) for every data row, and give
# the G2 header the same treatment as the other header cells.
# ---------------------------------------------------------------------
$ws.Range("E2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G2").Value = "Total"

$ws.Range("E3").Copy()
$ws.Range("G3:G26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($r = 3; $r -le 26; $r++) {
    $ws.Cells.Item($r, 7).Formula = "=E$r*D$r"
}

# ---------------------------------------------------------------------
# Row 26 previously held the "Total" formula (E26/G26); move that total
# down to a fresh row 32, now summing through the new rows.
# ---------------------------------------------------------------------
$ws.Range("E3").Copy()
$ws.Range("E32").PasteSpecial(-4122)
$ws.Range("G32").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E32").Value = "Total"
$ws.Range("G32").Formula = "=SUM(G3:G26)"

# ---------------------------------------------------------------------
# Row 20: note about the right-angle JST-PH connector used on this line.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = "Right angle JST-PH"

# ---------------------------------------------------------------------
# Selection cursor as left by the editing session.
# ---------------------------------------------------------------------
$ws.Range("C20").Select()
